$d = $word.ActiveDocument

# 1. Replace the lead-in text of the "Second local table" paragraph with the
#    new, longer text (keeps the trailing hyperlink + its formatting intact).
$d.Content.Find.Execute(
    "Second local table file: ",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Maecenas auctor lacus mauris, vitae fringilla nunc consequat eu. Interdum et malesuada fames ac ante ipsum primis in faucibus. Suspendisse faucibus, dolor vitae ultrices imperdiet, dui erat semper nisl, sed suscipit tellus risus a elit. This is the second local table link split in two lines: ",
    2)

# 2. Locate that paragraph again (now longer) and the two empty paragraphs
#    that trail it (the second of which carries the page-break run), then
#    remove the two empty paragraph marks so the page-break run becomes
#    part of the same paragraph as the hyperlink, exactly as in the target.
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*This is the data of the second table*") {
        $target = $i
        break
    }
}

$pNext = $d.Paragraphs.Item($target + 1)
$pNextNext = $d.Paragraphs.Item($target + 2)

# Delete the paragraph mark of the (empty) paragraph right after the
# hyperlink paragraph - merges it forward into the page-break paragraph.
$markA = $d.Range($pNext.Range.Start, $pNext.Range.End)
$markA.Delete()

# Delete the paragraph mark that used to end the hyperlink paragraph -
# merges the (now combined) page-break paragraph back into it.
$pTarget = $d.Paragraphs.Item($target)
$markB = $d.Range($pTarget.Range.End - 1, $pTarget.Range.End)
$markB.Delete()
